$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.04678466666666666
$ws.Range("H2").Value = 0.140354
$ws.Range("I2").Value = 0.006739448717762189
$ws.Range("J2").Value = 0.006739448717762188
$ws.Range("M2").Value = 3.626135
$ws.Range("N2").Value = 10.878405
$ws.Range("O2").Value = 0.4728835835086186
$ws.Range("P2").Value = 0.4728835835086186
$ws.Range("Q2").Value = 0.1696475172633333
$ws.Range("R2").Value = 1.52682765537
$ws.Range("S2").Value = 0.003186974660527948
$ws.Range("T2").Value = 0.003186974660527948

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.04678466666666666
$ws.Range("H3").Value = 0.140354
$ws.Range("I3").Value = 0.006739448717762189
$ws.Range("J3").Value = 0.006739448717762188
$ws.Range("O3").Value = 0.01581792773244636
$ws.Range("P3").Value = 0.01581792773244636
$ws.Range("Q3").Value = 0.005674699358666666
$ws.Range("R3").Value = 0.05107229422799999
$ws.Range("S3").Value = 0.0001066041127740906
$ws.Range("T3").Value = 0.0001066041127740906

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.04678466666666666
$ws.Range("H4").Value = 0.140354
$ws.Range("I4").Value = 0.006739448717762189
$ws.Range("J4").Value = 0.006739448717762188
$ws.Range("O4").Value = 0.511298488758935
$ws.Range("P4").Value = 0.511298488758935
$ws.Range("Q4").Value = 0.1834289077131111
$ws.Range("R4").Value = 1.650860169418
$ws.Range("S4").Value = 0.003445869944460149
$ws.Range("T4").Value = 0.003445869944460149

$ws.Range("I5").Value = 0.9176013393810419
$ws.Range("J5").Value = 0.9176013393810418
$ws.Range("M5").Value = 3.626135
$ws.Range("N5").Value = 10.878405
$ws.Range("O5").Value = 0.4728835835086186
$ws.Range("P5").Value = 0.4728835835086186
$ws.Range("Q5").Value = 23.09814876300333
$ws.Range("R5").Value = 207.88333886703
$ws.Range("S5").Value = 0.4339186095988152
$ws.Range("T5").Value = 0.4339186095988152

$ws.Range("I6").Value = 0.9176013393810419
$ws.Range("J6").Value = 0.9176013393810418
$ws.Range("O6").Value = 0.01581792773244636
$ws.Range("P6").Value = 0.01581792773244636
$ws.Range("S6").Value = 0.01451455167352531
$ws.Range("T6").Value = 0.01451455167352531

$ws.Range("I7").Value = 0.9176013393810419
$ws.Range("J7").Value = 0.9176013393810418
$ws.Range("O7").Value = 0.511298488758935
$ws.Range("P7").Value = 0.511298488758935
$ws.Range("S7").Value = 0.4691681781087014
$ws.Range("T7").Value = 0.4691681781087013

$ws.Range("G8").Value = 0.5252196666666668
$ws.Range("I8").Value = 0.07565921190119594
$ws.Range("J8").Value = 0.07565921190119593
$ws.Range("M8").Value = 3.626135
$ws.Range("N8").Value = 10.878405
$ws.Range("O8").Value = 0.4728835835086186
$ws.Range("P8").Value = 0.4728835835086186
$ws.Range("Q8").Value = 1.904517415988334
$ws.Range("R8").Value = 17.140656743895
$ws.Range("S8").Value = 0.03577799924927546
$ws.Range("T8").Value = 0.03577799924927546

$ws.Range("G9").Value = 0.5252196666666668
$ws.Range("I9").Value = 0.07565921190119594
$ws.Range("J9").Value = 0.07565921190119593
$ws.Range("O9").Value = 0.01581792773244636
$ws.Range("P9").Value = 0.01581792773244636
$ws.Range("Q9").Value = 0.06370599424866667
$ws.Range("R9").Value = 0.573353948238
$ws.Range("S9").Value = 0.001196771946146963
$ws.Range("T9").Value = 0.001196771946146963

$ws.Range("G10").Value = 0.5252196666666668
$ws.Range("I10").Value = 0.07565921190119594
$ws.Range("J10").Value = 0.07565921190119593
$ws.Range("O10").Value = 0.511298488758935
$ws.Range("P10").Value = 0.511298488758935
$ws.Range("Q10").Value = 2.059231723344778
$ws.Range("S10").Value = 0.03868444070577352
$ws.Range("T10").Value = 0.03868444070577351
